$d = $word.ActiveDocument

# --- Change 1: team-members line text fix, plus the "_GoBack" bookmark now
#     lives here (it tracks the most-recent edit location). ---
$teamRange = $d.Content
$found1 = $teamRange.Find.Execute(
    "חנין חטיב , פארוק, סגא, מחמוד", $true, $false, $false, $false, $false,
    $true, 1, $false, "חנין חטיב, פארוק, סגא", 2)
if (-not $found1) {
    throw "Could not find team-members paragraph text"
}

# --- Change 3 (handled before re-adding the bookmark): the Frontend tech
#     line had its text accidentally split by the old "_GoBack" bookmark
#     ("Materia" | bookmark | "l UI / Tailwind CSS]"). Remove the stray
#     bookmark and merge the text back into one run. ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$found3 = $d.Content.Find.Execute(
    " React, [Materia" + "l UI / Tailwind CSS]", $true, $false, $false,
    $false, $false, $true, 1, $false,
    " React, [Material UI / Tailwind CSS]", 2)
if (-not $found3) {
    throw "Could not find/fix the Frontend tech-stack line"
}

# Re-add "_GoBack" right after the team-members text we just edited (this is
# where Word now considers the last edit to have happened).
$teamRange.Find.Execute("חנין חטיב, פארוק, סגא") | Out-Null
$teamRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $teamRange) | Out-Null

# --- Change 2: drop the stray <w:rFonts w:hint="cs"/> from the paragraph
#     mark formatting of the (empty-looking) CHATBOT paragraph. There is no
#     direct object-model property for the raw "hint" attribute, so the
#     paragraph is rewritten via WordOpenXML/InsertXML with that one
#     attribute removed, keeping everything else byte-identical. ---
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("CHATBOT") -and $t.Contains("טקסטואלי")) {
        $target = $p
        break
    }
}
if ($null -eq $target) {
    throw "Could not find the CHATBOT paragraph"
}

$fixedXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="76D496EF" w14:textId="05533EC0" w:rsidR="00217B00" w:rsidRPr="00217B00" w:rsidRDefault="00217B00" w:rsidP="00217B00"><w:pPr><w:rPr><w:rtl/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="00217B00"><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:b/><w:bCs/><w:lang w:bidi="ar-AE"/></w:rPr><w:t xml:space="preserve">CHATBOT </w:t></w:r><w:r w:rsidRPr="00217B00"><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:b/><w:bCs/><w:rtl/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> טקסטואלי:</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> לתת אפשרות למשתמש לכתוב/ללחוץ על טקסט כדי להשתמש במערכת.</w:t></w:r></w:p>'
$target.Range.InsertXML($fixedXml)

Write-Host "Edit complete"
